$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to append (rows 17 to 24)
$data = @(
    @{ row=17; A=20200107; B=14; C=5; D=3; K=-1055; L=1283.82; M=841.79; N=45.3 },
    @{ row=18; A=20200107; B=14; C=5; D=3; K=-1054; L=1285.23; M=921.77; N=45 },
    @{ row=19; A=20200107; B=14; C=5; D=3; K=-1053; L=1284.5; M=1001.86; N=45.01 },
    @{ row=20; A=20200107; B=14; C=5; D=3; E=2020; F=3; G=10; H=18; I=14; J=49.562; K=-1002; L=1280; M=1008 },
    @{ row=21; A=20200107; B=14; C=5; D=3; E=2020; F=3; G=10; H=18; I=14; J=49.562; K=-1001; L=1280; M=985 },
    @{ row=22; A=20200107; B=14; C=5; D=3; E=2020; F=3; G=10; H=18; I=14; J=49.562; K=-1000; L=1280; M=963 },
    @{ row=23; A=20200107; B=14; C=5; D=3; E=2020; F=3; G=10; H=18; I=14; J=49.562; K=-999; L=1283; M=942 },
    @{ row=24; A=20200107; B=14; C=5; D=3; E=2020; F=3; G=10; H=18; I=14; J=49.562; K=-998; L=1284; M=922 }
)

foreach ($r in $data) {
    $row = $r.row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    if ($r.ContainsKey('E')) { $ws.Cells.Item($row, 5).Value = $r.E }
    if ($r.ContainsKey('F')) { $ws.Cells.Item($row, 6).Value = $r.F }
    if ($r.ContainsKey('G')) { $ws.Cells.Item($row, 7).Value = $r.G }
    if ($r.ContainsKey('H')) { $ws.Cells.Item($row, 8).Value = $r.H }
    if ($r.ContainsKey('I')) { $ws.Cells.Item($row, 9).Value = $r.I }
    if ($r.ContainsKey('J')) {
        $cell = $ws.Cells.Item($row, 10)
        $cell.Value = $r.J
        $cell.NumberFormat = "0.00E+00"
    }
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    if ($r.ContainsKey('N')) { $ws.Cells.Item($row, 14).Value = $r.N }
}

# Update the view: scroll so A4 is the top-left visible cell, and select Q23
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("Q23").Select()
